$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays text-formatted so values like "1.005" are not
# reinterpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.628.19"
$ws.Range("E2").Value = "  -7.37%  "

$ws.Range("D3").Value = "1.697.10"
$ws.Range("E3").Value = "  -6.02%  "

$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").Value = "219.46"
$ws.Range("E5").Value = "  -5.54%  "

$ws.Range("D6").Value = "0.5128"
$ws.Range("E6").Value = "  -13.16%  "

$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "0.2657"
$ws.Range("E8").Value = "  -4.40%  "

$ws.Range("D9").Value = "22.16"
$ws.Range("E9").Value = "  -4.87%  "

$ws.Range("D10").Value = "0.06252"
$ws.Range("E10").Value = "  -8.46%  "

$ws.Range("D11").Value = "0.07323"
$ws.Range("E11").Value = "  -2.40%  "

$ws.Range("D12").Value = "1.698.68"
$ws.Range("E12").Value = "  -5.98%  "

$ws.Range("D13").Value = "4.511"
$ws.Range("E13").Value = "  -5.31%  "

$ws.Range("D14").Value = "0.5839"
$ws.Range("E14").Value = "  -6.30%  "

$ws.Range("D15").Value = "1.930.38"
$ws.Range("E15").Value = "  -5.89%  "

$ws.Range("D16").Value = "0.000008410"
$ws.Range("E16").Value = "  -9.37%  "

$ws.Range("D17").Value = "65.52"
$ws.Range("E17").Value = "  -13.45%  "

$ws.Range("D18").Value = "26.685.61"
$ws.Range("E18").Value = "  -7.04%  "

$ws.Range("D19").Value = "5.051"
$ws.Range("E19").Value = "  -7.77%  "

$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("D21").Value = "10.88"
$ws.Range("E21").Value = "  -5.23%  "

$ws.Range("D22").Value = "187.31"
$ws.Range("E22").Value = "  -11.40%  "

$ws.Range("D23").Value = "6.270"
$ws.Range("E23").Value = "  -8.33%  "

$ws.Range("D24").Value = "1.006"
$ws.Range("E24").Value = "  +0.20%  "

$ws.Range("D25").Value = "145.19"
$ws.Range("E25").Value = "  -5.86%  "

$ws.Range("D26").Value = "7.615"
$ws.Range("E26").Value = "  -3.28%  "

$ws.Range("D27").Value = "0.1150"

$ws.Range("D28").Value = "15.79"
$ws.Range("E28").Value = "  -3.95%  "

$ws.Range("E29").Value = "  -9.04%  "

$ws.Range("D30").Value = "0.05717"
$ws.Range("E30").Value = "  -7.36%  "

$ws.Range("E31").Value = "  -6.18%  "

$ws.Range("D32").Value = "3.521"
$ws.Range("E32").Value = "  -6.27%  "

$ws.Range("D33").Value = "3.507"
$ws.Range("E33").Value = "  -7.34%  "

$ws.Range("D34").Value = "1.665"
$ws.Range("E34").Value = "  -4.07%  "

$ws.Range("D35").Value = "1.023"
$ws.Range("E35").Value = "  -3.92%  "

$ws.Range("D36").Value = "0.6017"
$ws.Range("E36").Value = "  -6.39%  "

$ws.Range("D37").Value = "2.374"
$ws.Range("E37").Value = "  -4.86%  "

$ws.Range("D38").Value = "2.684"
$ws.Range("E38").Value = "  -1.22%  "

$ws.Range("D39").Value = "1.095.29"
$ws.Range("E39").Value = "  -4.34%  "

$ws.Range("D40").Value = "0.01600"
$ws.Range("E40").Value = "  -5.97%  "

$ws.Range("D41").Value = "0.8654"
$ws.Range("E41").Value = "  -1.99%  "

$ws.Range("D42").Value = "5.902"
$ws.Range("E42").Value = "  -10.04%  "

$ws.Range("D43").Value = "1.003"
$ws.Range("E43").Value = "  -0.35%  "

$ws.Range("D44").Value = "98.69"
$ws.Range("E44").Value = "  -1.43%  "

$ws.Range("D45").Value = "1.858.82"

$ws.Range("D46").Value = "0.00000000109"
$ws.Range("E46").Value = "  -2.71%  "

$ws.Range("D47").Value = "56.82"
$ws.Range("E47").Value = "  -6.10%  "

$ws.Range("D48").Value = "8.221"
$ws.Range("E48").Value = "  -1.77%  "

$ws.Range("D49").Value = "1.004"
$ws.Range("E49").Value = "  -0.27%  "

$ws.Range("D50").Value = "0.05244"
$ws.Range("E50").Value = "  -4.17%  "

$ws.Range("D51").Value = "0.4324"
$ws.Range("E51").Value = "  -3.62%  "
